$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the C1:C3 cells (same style used across the data rows)
# onto the new D1:D3 cells so the added "Result" labels match the existing look.
$ws.Range("C1:C3").Copy() | Out-Null
$ws.Range("D1:D3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Add the new "Result" column values for the first three data rows.
$ws.Range("D1").Value = "Result1"
$ws.Range("D2").Value = "Result2"
$ws.Range("D3").Value = "Result3"
